$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-140 are updated from serial date 45189 (2023-09-20)
# to serial date 45190 (2023-09-21).
$ws.Range("C2:C140").Value = 45190
